{"js": "const pairs = [\n  [\n    \"Hello again! Kuba ngumtali kungaba matima, kungako kufanele utinakekele!\",\n    \"Sawubona futsi! Kuba ngumtali kungaba matima, kungako kufanele utinakekele!\"\n  ],\n  [\n    \"Here is a simple stretching and movement activity that may help you with stress.\",\n    \"Nayi indlela lemelula yekwelula nekuvocavoca umtimba lengakusita kwehlisa kucindzeteleka.\"\n  ],\n  [\n    \"Stand up and stretch your arms up to the sky.\",\n    \"Sukuma welule imikhono yakho netandla kube shangatsi ufuna kutsintsa sibhakabhaka.\"\n  ],\n  [\n    \"Push away any clouds and let the sun shine down on you.\",\n    \"Yenta shangatsi ukhweshisa emafu uvulela umsebe welilanga kutsi ukukhanyise.\"\n  ],\n  [\n    \"Stretch to both sides.\",\n    \"Yelula umtimba ubhekise kuto totimbili tinhlangotsi, ngesekuncele nangesekudla.\"\n  ],\n  [\n    \"Stretch to the front and to the back.\",\n    \"Yelula ubhekise ngaphambili nangemuva.\"\n  ],\n  [\n    \"Let your hands hang by your side.\",\n    \"Yekela tandla tilengele emaceleni.\"\n  ],\n  [\n    \"Squeeze your shoulders tightly up to your ears as you breathe in, scrunch your eyes, and hold your breath.\",\n    \"Cindzetela emahlombe akho kakhulu uwaphakamise ate afike etindlebeni ube udvonsa umoya, uswace, uchubeke nekubamba umoya ungawukhiphi.\"\n  ],\n  [\n    \"Release your shoulders as you relax your body and breathe out.\",\n    \"Yehlisa emahlombe uphumute umtimba wakho bese ukhipha umoya.\"\n  ],\n  [\n    \"Relax your arms and let them swing to the front and back. Take 4 deep breaths while you do this. This is good for your lower back.\",\n    \"Jikitisa emahlombe akho uwayise embili nangemuva. Dvonsa umoya uphindze uwukhiphe kute kube mahlandla lamane uselula emahlombe. Loku kuyakusita kwelula lukhalo.\"\n  ],\n  [\n    \"Relax your arms and let your arms swing sideways, and turn your upper body. Your arms should gently hit your back as you twist from side to side.\",\n    \"Yehlisa emahlombe bese uyawajikitisa uwabhekise emaceleni ngasemhlubulweni nasetulu. Imikhono yakho kufanele iwutsintse kancane nje umhlane ngesikhatsi uyijikitisela emaceleni.\"\n  ],\n  [\n    \"Place one foot in front of the other and make small circles from the ankles, remember to turn to both sides. Do each foot with four circles in each direction.\",\n    \"Beka lunyawo lunye phambi kwalolunye, unyatsele ngetintwane uphakamise sitsendze, uhambise lunyawo wente indingilizi lencane. Khumbula kushintja lunyawo. Yenta lendingilizi kute kube kane ubhekise ngesekudla bese ubhekisa ngesekuncele.\"\n  ],\n  [\n    \"Now shake your whole body as fast as you can. Shake it up high! Shake it down low! Shake it all around!\",\n    \"Nyalo shukumisa umtimba wakho wonkhe usheshise kakhulu. Wushukumise uye etulu! Wushukumise ushone phansi kakhulu! Wushukumise ubheke yonkhe indzawo!\"\n  ],\n  [\n    \"Raise your arms slowly above your head as you breathe in. Lower your arms slowly to your side as you breathe out.\",\n    \"Phakamisa imikhono kancane kancane ibe ngetulu kwenhloko ube udvonsa umoya. Yehlisa imikhono kancane kancane ube ukhiphe umoya.\"\n  ],\n  [\n    \"Raise and lower your arms, breathing in and out,  3 more times.\",\n    \"Phakamisa uphindze wehlise imikhono, udvonse uphindze ukhiphe umoya, kute kube katsatfu.\"\n  ],\n  [\n    \"Take a moment to reflect on your experience.\",\n    \"Tsatsa sikhashana ucabanga ngaloku locedza kukwenta.\"\n  ],\n  [\n    \"Remember that you can do this activity at any time whenever you feel like you need to release stress and energise your mind and body.\",\n    \"Khumbula kutsi ungakwenta loku noma ngabe ngunini nawutiva kutsi udzinga kwehlisa kucindzeteleka phindze uphaphamise umtimba kanye nengcondvo.\"\n  ],\n  [\n    \"Your home activity is to try to try to do this activity at least once every day. Ungayetama nyalo?\",\n    \"Umsebenti wasekhaya kutsi wetame kuvocavoca umtimba lokungenani kanye ngelilanga. Ungayetama nyalo?\"\n  ],\n  [\n    \"Home Activity: Stretch and move every day. \",\n    \"Umsebenti Wasekhaya: Telule uphindze utivocavoce njalo ngelilanga. \"\n  ],\n  [\n    \"Module: Mental Changes in the Teen Years\",\n    \"Module: Kushintja Kwengcondvo Eminyakeni Yekutfomba\"\n  ],\n  [\n    \"Module: Social Changes in the Teen Years\",\n    \"Kushintja Kwemphilo Eminyakeni Yekutfomba\"\n  ]\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Each source string is expected to appear exactly once; replace every\n  // match found (defensively) by setting its text in place.\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"Hello again! Kuba ngumtali kungaba matima, kungako kufanele utinakekele!\", \"Sawubona futsi! Kuba ngumtali kungaba matima, kungako kufanele utinakekele!\")\n  ,@(\"Here is a simple stretching and movement activity that may help you with stress.\", \"Nayi indlela lemelula yekwelula nekuvocavoca umtimba lengakusita kwehlisa kucindzeteleka.\")\n  ,@(\"Stand up and stretch your arms up to the sky.\", \"Sukuma welule imikhono yakho netandla kube shangatsi ufuna kutsintsa sibhakabhaka.\")\n  ,@(\"Push away any clouds and let the sun shine down on you.\", \"Yenta shangatsi ukhweshisa emafu uvulela umsebe welilanga kutsi ukukhanyise.\")\n  ,@(\"Stretch to both sides.\", \"Yelula umtimba ubhekise kuto totimbili tinhlangotsi, ngesekuncele nangesekudla.\")\n  ,@(\"Stretch to the front and to the back.\", \"Yelula ubhekise ngaphambili nangemuva.\")\n  ,@(\"Let your hands hang by your side.\", \"Yekela tandla tilengele emaceleni.\")\n  ,@(\"Squeeze your shoulders tightly up to your ears as you breathe in, scrunch your eyes, and hold your breath.\", \"Cindzetela emahlombe akho kakhulu uwaphakamise ate afike etindlebeni ube udvonsa umoya, uswace, uchubeke nekubamba umoya ungawukhiphi.\")\n  ,@(\"Release your shoulders as you relax your body and breathe out.\", \"Yehlisa emahlombe uphumute umtimba wakho bese ukhipha umoya.\")\n  ,@(\"Relax your arms and let them swing to the front and back. Take 4 deep breaths while you do this. This is good for your lower back.\", \"Jikitisa emahlombe akho uwayise embili nangemuva. Dvonsa umoya uphindze uwukhiphe kute kube mahlandla lamane uselula emahlombe. Loku kuyakusita kwelula lukhalo.\")\n  ,@(\"Relax your arms and let your arms swing sideways, and turn your upper body. Your arms should gently hit your back as you twist from side to side.\", \"Yehlisa emahlombe bese uyawajikitisa uwabhekise emaceleni ngasemhlubulweni nasetulu. Imikhono yakho kufanele iwutsintse kancane nje umhlane ngesikhatsi uyijikitisela emaceleni.\")\n  ,@(\"Place one foot in front of the other and make small circles from the ankles, remember to turn to both sides. Do each foot with four circles in each direction.\", \"Beka lunyawo lunye phambi kwalolunye, unyatsele ngetintwane uphakamise sitsendze, uhambise lunyawo wente indingilizi lencane. Khumbula kushintja lunyawo. Yenta lendingilizi kute kube kane ubhekise ngesekudla bese ubhekisa ngesekuncele.\")\n  ,@(\"Now shake your whole body as fast as you can. Shake it up high! Shake it down low! Shake it all around!\", \"Nyalo shukumisa umtimba wakho wonkhe usheshise kakhulu. Wushukumise uye etulu! Wushukumise ushone phansi kakhulu! Wushukumise ubheke yonkhe indzawo!\")\n  ,@(\"Raise your arms slowly above your head as you breathe in. Lower your arms slowly to your side as you breathe out.\", \"Phakamisa imikhono kancane kancane ibe ngetulu kwenhloko ube udvonsa umoya. Yehlisa imikhono kancane kancane ube ukhiphe umoya.\")\n  ,@(\"Raise and lower your arms, breathing in and out,  3 more times.\", \"Phakamisa uphindze wehlise imikhono, udvonse uphindze ukhiphe umoya, kute kube katsatfu.\")\n  ,@(\"Take a moment to reflect on your experience.\", \"Tsatsa sikhashana ucabanga ngaloku locedza kukwenta.\")\n  ,@(\"Remember that you can do this activity at any time whenever you feel like you need to release stress and energise your mind and body.\", \"Khumbula kutsi ungakwenta loku noma ngabe ngunini nawutiva kutsi udzinga kwehlisa kucindzeteleka phindze uphaphamise umtimba kanye nengcondvo.\")\n  ,@(\"Your home activity is to try to try to do this activity at least once every day. Ungayetama nyalo?\", \"Umsebenti wasekhaya kutsi wetame kuvocavoca umtimba lokungenani kanye ngelilanga. Ungayetama nyalo?\")\n  ,@(\"Home Activity: Stretch and move every day. \", \"Umsebenti Wasekhaya: Telule uphindze utivocavoce njalo ngelilanga. \")\n  ,@(\"Module: Mental Changes in the Teen Years\", \"Module: Kushintja Kwengcondvo Eminyakeni Yekutfomba\")\n  ,@(\"Module: Social Changes in the Teen Years\", \"Kushintja Kwemphilo Eminyakeni Yekutfomba\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n"}
